# Generate Report for Archive
# The handback/status report is being regenerated: every cell that still
# shows the old "Ready for handoff" status is moved along to "In Translation",
# and the (now shorter) Status columns are re-sized to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E & F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsOverview.Columns("E:F").ColumnWidth = 13.4101845877511

# --- zh-cn sheet: Status column (C) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C4").Value = "In Translation"
$wsZh.Columns("C:C").ColumnWidth = 13.4101845877511

# --- de-de sheet: Status column (C) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C4").Value = "In Translation"
$wsDe.Columns("C:C").ColumnWidth = 13.4101845877511
